$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2350.625
$ws.Range("I70").Value = 2700.8333
$ws.Range("J70").Value = 1300
$ws.Range("K70").Value = 8102.499899999999
$ws.Range("L70").Value = 3900
$ws.Range("M70").Value = -7832.499899999999
$ws.Range("N70").Value = -4440

$ws.Range("H73").Value = 2350.625
$ws.Range("I73").Value = 2700.8333
$ws.Range("J73").Value = 1300
$ws.Range("K73").Value = 8102.499899999999
$ws.Range("L73").Value = 3900
$ws.Range("M73").Value = -7166.499899999999
$ws.Range("N73").Value = -5772

$ws.Range("H131").Value = 1918.4348
$ws.Range("I131").Value = 554.94116
$ws.Range("J131").Value = 5781.6665
$ws.Range("K131").Value = 1664.82348
$ws.Range("L131").Value = 17344.9995
$ws.Range("M131").Value = 3375.17652
$ws.Range("N131").Value = -27424.9995

$ws.Range("H137").Value = 15152828
$ws.Range("I137").Value = 25001044
$ws.Range("J137").Value = 1726.6154
$ws.Range("K137").Value = 75003132
$ws.Range("L137").Value = 5179.8462
$ws.Range("M137").Value = -75000582
$ws.Range("N137").Value = -10279.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 5743.75
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5743.75
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 5743.75
$ws.Range("N43").Value = -6369.75

$ws.Range("H61").Value = 6174160
$ws.Range("I61").Value = 6668013
$ws.Range("J61").Value = 1000
$ws.Range("K61").Value = 6668013
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = -6667801
$ws.Range("N61").Value = -1424

$ws.Range("H122").Value = 1418.375
$ws.Range("I122").Value = 1363.8572
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 4091.5716
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -1641.5716
$ws.Range("N122").Value = -10300

$ws.Range("H136").Value = 6174160
$ws.Range("I136").Value = 6668013
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 20004039
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -20001489
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 111113110
$ws.Range("I105").Value = 2336.3333
$ws.Range("J105").Value = 166668500
$ws.Range("K105").Value = 2336.3333
$ws.Range("L105").Value = 166668500
$ws.Range("M105").Value = -589.3332999999998
$ws.Range("N105").Value = -166671994

$ws.Range("H134").Value = 2808.2322
$ws.Range("I134").Value = 922.23914
$ws.Range("J134").Value = 11483.8
$ws.Range("K134").Value = 2766.71742
$ws.Range("L134").Value = 34451.39999999999
$ws.Range("M134").Value = -231.7174199999999
$ws.Range("N134").Value = -39521.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2428.4285
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 2666.5
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 7999.5
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -9371.5

$ws.Range("H65").Value = 2428.4285
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 2666.5
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 23998.5
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -30862.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 42200
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 42200
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 42200
$ws.Range("N103").Value = -44544

$ws.Range("H123").Value = 10296.363
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10296.363
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10296.363
$ws.Range("N123").Value = -15196.363

$ws.Range("H126").Value = 3460
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 3950
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 11850
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -16790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 41667916
$ws.Range("I40").Value = 1500.8
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 1500.8
$ws.Range("L40").Value = 250000000
$ws.Range("M40").Value = -1364.8
$ws.Range("N40").Value = -250000272

$ws.Range("H122").Value = 22381.8
$ws.Range("I122").Value = 26726
$ws.Range("J122").Value = 5005
$ws.Range("K122").Value = 80178
$ws.Range("L122").Value = 15015
$ws.Range("M122").Value = -77728
$ws.Range("N122").Value = -19915

$ws.Range("H129").Value = 59800
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 59800
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 59800
$ws.Range("N129").Value = -69800

$ws.Range("H132").Value = 35724284
$ws.Range("I132").Value = 58825904
$ws.Range("J132").Value = 21782.182
$ws.Range("K132").Value = 176477712
$ws.Range("L132").Value = 65346.546
$ws.Range("M132").Value = -176475182
$ws.Range("N132").Value = -70406.546

$ws.Range("H136").Value = 35715970
$ws.Range("I136").Value = 6804738
$ws.Range("J136").Value = 90910136
$ws.Range("K136").Value = 20414214
$ws.Range("L136").Value = 272730408
$ws.Range("M136").Value = -20411664
$ws.Range("N136").Value = -272735508

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 38200
$ws.Range("I82").Value = 19800
$ws.Range("J82").Value = 47400
$ws.Range("K82").Value = 19800
$ws.Range("L82").Value = 47400
$ws.Range("M82").Value = -19417
$ws.Range("N82").Value = -48166

$ws.Range("H85").Value = 38200
$ws.Range("I85").Value = 19800
$ws.Range("J85").Value = 47400
$ws.Range("K85").Value = 19800
$ws.Range("L85").Value = 47400
$ws.Range("M85").Value = -18474
$ws.Range("N85").Value = -50052

$ws.Range("H122").Value = 26062.809
$ws.Range("I122").Value = 43001.5
$ws.Range("J122").Value = 3477.889
$ws.Range("K122").Value = 129004.5
$ws.Range("L122").Value = 10433.667
$ws.Range("M122").Value = -126554.5
$ws.Range("N122").Value = -15333.667

$ws.Range("H126").Value = 1229.4445
$ws.Range("I126").Value = 717
$ws.Range("J126").Value = 1870
$ws.Range("K126").Value = 2151
$ws.Range("L126").Value = 5610
$ws.Range("M126").Value = 319
$ws.Range("N126").Value = -10550

$ws.Range("H132").Value = 12224650
$ws.Range("I132").Value = 32988.938
$ws.Range("J132").Value = 62515250
$ws.Range("K132").Value = 98966.81400000001
$ws.Range("L132").Value = 187545750
$ws.Range("M132").Value = -96436.81400000001
$ws.Range("N132").Value = -187550810

$ws.Range("H136").Value = 20002182
$ws.Range("I136").Value = 27779114
$ws.Range("J136").Value = 4357.143
$ws.Range("K136").Value = 83337342
$ws.Range("L136").Value = 13071.429
$ws.Range("M136").Value = -83334792
$ws.Range("N136").Value = -18171.429
